$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "21.83")
# are preserved exactly as text instead of being parsed into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.013.38"
$ws.Range("E2").Value = "  -2.24%  "
$ws.Range("D3").Value = "1.668.21"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "216.85"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "0.2656"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.06410"
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").Value = "21.83"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D11").Value = "0.07437"
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("D12").Value = "1.670.45"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "0.5841"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "0.000008575"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "64.42"
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("D17").Value = "26.077.29"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("D21").Value = "192.47"
$ws.Range("E21").Value = "  +3.34%  "
$ws.Range("D22").Value = "6.209"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").Value = "1.006"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "144.74"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "7.622"
$ws.Range("E25").Value = "  +2.05%  "
$ws.Range("D26").Value = "0.1197"
$ws.Range("E26").Value = "  +2.87%  "
$ws.Range("E27").Value = "  -1.31%  "
$ws.Range("D28").Value = "0.06502"
$ws.Range("E28").Value = "  +13.80%  "
$ws.Range("D29").Value = "1.327"
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("D31").Value = "3.544"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("D33").Value = "1.650"
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("D34").Value = "1.019"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "0.6107"
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("D36").Value = "2.367"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("D37").Value = "2.690"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "6.257"
$ws.Range("E38").Value = "  +7.25%  "
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").Value = "1.092.04"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").Value = "0.8632"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("D43").Value = "100.79"
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").Value = "1.817.47"
$ws.Range("E44").Value = "  -1.84%  "
$ws.Range("E45").Value = "  -1.64%  "
$ws.Range("D46").Value = "56.42"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").Value = "1.010"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("D48").Value = "8.088"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("D49").Value = "0.05237"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").Value = "6.041"
$ws.Range("E51").Value = "  +4.13%  "
